$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.251
$ws.Range("B12").Value = 5.492999999999999
$ws.Range("E14").Value = 17.165
$ws.Range("E26").Value = 16.117
$ws.Range("E31").Value = 16.232
$ws.Range("B32").Value = 6.453
$ws.Range("E35").Value = 16.63
$ws.Range("B36").Value = 7.987
$ws.Range("E37").Value = 16.737
$ws.Range("B38").Value = 5.342000000000001
$ws.Range("E45").Value = 16.749
$ws.Range("B46").Value = 6.334
$ws.Range("B54").Value = 5.155
$ws.Range("B55").Value = 4.683999999999999
$ws.Range("E57").Value = 16.529
$ws.Range("B67").Value = 5.286
$ws.Range("B69").Value = 5.002000000000001
$ws.Range("B72").Value = 5.380999999999999
$ws.Range("B91").Value = 6.219
$ws.Range("B99").Value = 5.217
$ws.Range("E100").Value = 16.749
$ws.Range("E102").Value = 16.436
